# Delete entity-management test case: "duplicated additional properties"
# (row that was createEntities test "Test-15" data, i.e. physical row 16).
# Columns B:I of the subsequent rows (17-21) shift up into rows (16-20),
# while column A (the sequential Test-N label) is left untouched. The
# now-empty trailing row 21 is then removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("createEntities")

$cols = @("B","C","D","E","F","G","H","I")

for ($destRow = 16; $destRow -le 20; $destRow++) {
    $srcRow = $destRow + 1
    foreach ($col in $cols) {
        $srcCell = $ws.Range("$col$srcRow")
        $dstCell = $ws.Range("$col$destRow")
        $dstCell.Value = $srcCell.Value()
    }
}

# Remove the now-duplicated last row (21), which has shifted up already.
$ws.Range("A21:I21").EntireRow.Delete()

# Restore focus to the createEntities sheet / tab as the active sheet.
$ws.Select()
$ws.Range("A2:A20").Select()
